$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '63.271.31'
    'E2' = '  +6.37%  '
    'D3' = '3.111.70'
    'E3' = '  +4.29%  '
    'E4' = '  +0.06%  '
    'D5' = '585.18'
    'E5' = '  +3.46%  '
    'D6' = '144.95'
    'E6' = '  +4.60%  '
    'E7' = '  +0.05%  '
    'D8' = '3.102.80'
    'E8' = '  +4.28%  '
    'D9' = '0.530'
    'E9' = '  +1.82%  '
    'E10' = '  +13.47%  '
    'D11' = '5.78'
    'E11' = '  +7.39%  '
    'E12' = '  +4.09%  '
    'E13' = '  +7.88%  '
    'D14' = '35.59'
    'E14' = '  +5.47%  '
    'E15' = '  +0.64%  '
    'D16' = '3.630.95'
    'E16' = '  +4.51%  '
    'D17' = '7.17'
    'E17' = '  +1.70%  '
    'D18' = '63.195.51'
    'E18' = '  +6.25%  '
    'D19' = '3.110.02'
    'E19' = '  +4.39%  '
    'D20' = '467.20'
    'E20' = '  +7.02%  '
    'D21' = '14.15'
    'E21' = '  +3.87%  '
    'D22' = '0.727'
    'E22' = '  +0.80%  '
    'D23' = '7.53'
    'E23' = '  +7.27%  '
    'D24' = '13.32'
    'E24' = '  +0.07%  '
    'D25' = '82.11'
    'E26' = '  +0.02%  '
    'E27' = '  +10.41%  '
    'D28' = '2.22'
    'E28' = '  -0.22%  '
    'E29' = '  +4.96%  '
    'E30' = '  +0.03%  '
    'E31' = '  +10.95%  '
    'D32' = '26.94'
    'E32' = '  +4.63%  '
    'E33' = '  +5.09%  '
    'D34' = '0.0₃0865'
    'E34' = '  +11.55%  '
    'D35' = '2.41'
    'E35' = '  +16.40%  '
    'E36' = '  +7.11%  '
    'B37' = 'Filecoin'
    'C37' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D37' = '6.08'
    'E37' = '  +2.95%  '
    'B38' = 'dogwifhat'
    'C38' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D38' = '3.33'
    'E38' = '  +19.32%  '
    'D39' = '50.57'
    'E39' = '  +4.01%  '
    'D40' = '439.90'
    'E40' = '  +10.31%  '
    'D41' = '8.72'
    'E41' = '  +0.57%  '
    'D42' = '2.923.66'
    'E42' = '  +6.94%  '
    'E43' = '  +5.22%  '
    'E44' = '  +11.73%  '
    'E45' = '  +4.46%  '
    'E46' = '  +8.31%  '
    'B47' = 'USDe'
    'C47' = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
    'D47' = '0.999'
    'E47' = '  +0.02%  '
    'B48' = 'Arweave'
    'C48' = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
    'D48' = '35.04'
    'E48' = '  -0.54%  '
    'D49' = '123.10'
    'E49' = '  +0.79%  '
    'E50' = '  +0.94%  '
    'D51' = '24.69'
    'E51' = '  +6.09%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $updates[$cellRef]
    $cell.Style = "Normal"
}
